$d = $word.ActiveDocument

# 1. Merge "ASSIGNMENT " + "2" into a single run "ASSIGNMENT 2"
$d.Content.Find.Execute("ASSIGNMENT 2", $false, $false, $false, $false, $false, $true, 1, $false, "ASSIGNMENT 2", 2) | Out-Null
